$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we touch so that numeric-looking
# strings (e.g. '594.00', '0.110', '2.69') are preserved exactly as text
# instead of being auto-converted to numbers and losing formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.899.79'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.034.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.00'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.58'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.80%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.028.35'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.57'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +11.43%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.83%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.536.95'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.12'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.860.17'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.036.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '453.02'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.28'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +7.09%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.44'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.69'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.69%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.57'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0869'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.71%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.33%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +10.96%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.55'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.10'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +12.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.78'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '396.97'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.35%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.729.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.29'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.28'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.44'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.39%  '
